# Simulator full-month coverage / persist logs / fix employees
# Updates client names, employee id, and fills in rate/total numbers that
# were previously zeroed out (simulator now produces real hours -> $ data).

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Timesheet")
$wsSchema = $wb.Worksheets.Item("Jason Schema")

# ---- Weekly Timesheet: fix client names on the daily rows ----
$wsWeekly.Range("B2").Value = "McClure"
$wsWeekly.Range("B3").Value = "Evans"
$wsWeekly.Range("B4").Value = "Fritts"
$wsWeekly.Range("B5").Value = "Hendricks"
$wsWeekly.Range("B6").Value = "Regan"

# ---- Weekly Timesheet: fill in Rate / Total for each daily row ----
$wsWeekly.Range("E2").Value = 95
$wsWeekly.Range("F2").Value = 760
$wsWeekly.Range("E3").Value = 95
$wsWeekly.Range("F3").Value = 760
$wsWeekly.Range("E4").Value = 95
$wsWeekly.Range("F4").Value = 760
$wsWeekly.Range("E5").Value = 95
$wsWeekly.Range("F5").Value = 760
$wsWeekly.Range("E6").Value = 95
$wsWeekly.Range("F6").Value = 760

# ---- Weekly Timesheet: subtotal / grand total rows ----
$wsWeekly.Range("F8").Value = 3800
$wsWeekly.Range("F11").Value = 3800
$wsWeekly.Range("F13").Value = 3800

# ---- Jason Schema: same client-name fixes (one row per day) ----
$wsSchema.Range("D2").Value = "McClure"
$wsSchema.Range("D3").Value = "Evans"
$wsSchema.Range("D4").Value = "Fritts"
$wsSchema.Range("D5").Value = "Hendricks"
$wsSchema.Range("D6").Value = "Regan"

# ---- Jason Schema: fill in Rate / Total for each row ----
$wsSchema.Range("F2").Value = 95
$wsSchema.Range("G2").Value = 760
$wsSchema.Range("F3").Value = 95
$wsSchema.Range("G3").Value = 760
$wsSchema.Range("F4").Value = 95
$wsSchema.Range("G4").Value = 760
$wsSchema.Range("F5").Value = 95
$wsSchema.Range("G5").Value = 760
$wsSchema.Range("F6").Value = 95
$wsSchema.Range("G6").Value = 760

# ---- Jason Schema: fix employee id (shared across B2:B6) ----
$wsSchema.Range("B2:B6").Value = "emp_4nlnrvy7"
